# Add three new withdrawal rows (09 and 10 May 2018) to the bottom of the
# data, reusing the formatting of the last existing data row (row 15) so the
# new rows' number formats / styles match the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 15

for ($i = 0; $i -lt 3; $i++) {
    $ws.Rows($lastRow).Copy()
    $ws.Rows($lastRow + 1).Insert(-4121)  # xlShiftDown, formatting comes from copied row
    $lastRow = $lastRow + 1
}

# Row 16 : 09/05/2018 - Lanzi Patrizia - Tessuto n 10 - Mt. - 0.8
$ws.Cells.Item(16, 1).Value = "5/9/2018"
$ws.Cells.Item(16, 2).Value = "Lanzi Patrizia"
$ws.Cells.Item(16, 3).Value = "Tessuto n 10"
$ws.Cells.Item(16, 4).Value = "Mt."
$ws.Cells.Item(16, 5).Value = 0.8

# Row 17 : 09/05/2018 - Di Giacomo Caterina - Scalimetri - N°. - 43
$ws.Cells.Item(17, 1).Value = "5/9/2018"
$ws.Cells.Item(17, 2).Value = "Di Giacomo Caterina"
$ws.Cells.Item(17, 3).Value = "Scalimetri"
$ws.Cells.Item(17, 4).Value = "N°."
$ws.Cells.Item(17, 5).Value = 43

# Row 18 : 10/05/2018 - Puzziferri Domenico - Adesivo Pesante Bianco - Mt. - 1.5
$ws.Cells.Item(18, 1).Value = "5/10/2018"
$ws.Cells.Item(18, 2).Value = "Puzziferri Domenico"
$ws.Cells.Item(18, 3).Value = "Adesivo Pesante Bianco"
$ws.Cells.Item(18, 4).Value = "Mt."
$ws.Cells.Item(18, 5).Value = 1.5
